$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (title reflects "through" date)
$ws.Name = "Through 2022-07-22"

# Update the label for the July row (shared string)
$ws.Range("A8").Value = "July (through 07-22)"

# Update July row (row 8) data for columns C..I
$ws.Range("C8").Value = 44
$ws.Range("D8").Value = 48
$ws.Range("E8").Value = 57
$ws.Range("F8").Value = 35
$ws.Range("G8").Value = 94
$ws.Range("H8").Value = 107
$ws.Range("I8").Value = 125

# Update Total row (row 9) data for columns C..I
$ws.Range("C9").Value = 292
$ws.Range("D9").Value = 438
$ws.Range("E9").Value = 410
$ws.Range("F9").Value = 286
$ws.Range("G9").Value = 566
$ws.Range("H9").Value = 867
$ws.Range("I9").Value = 931
